$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated stat values (sval data regenerated to filter save games)
$data = @{
    2 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    3 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    4 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
    5 = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.896700893398075)
    6 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.034748368925986)
    7 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
